$wb = $excel.ActiveWorkbook

# Delete the "PA3" sheet entirely.
$wb.Worksheets.Item("PA3").Delete()

$ws = $wb.Worksheets.Item("PA")
$ws.Activate()

# Append a third overview-review bullet to the "PA4" column (column E) task
# descriptions for each student (rows 2-6).
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 + "`n· Đánh giá tổng quan lần 3"
}

# The extra line pushes the wrapped-text rows taller; reflect Excel's
# auto row-height recalculation for the affected rows.
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 75
